# Microcontroller Team 1.xlsx -- apply the authored edit:
#   1. Add a new row-5 label in column C describing the state diagram.
#   2. Move the active selection to F14 (where the user ended up after
#      adding/looking at the new content).
#   3. Best-effort: rename the built-in "Normal" cell style to "Standard"
#      (cosmetic artifact of the authoring app's locale).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label cell (becomes shared-string index 8 in sharedStrings.xml)
$ws.Range("C5").Value = "State Diagram for Pedestrian And Car TLS"

# Built-in "Normal" style -> "Standard" (matches the localized name Excel
# uses for the default style in some locales). Best effort: harmless if
# the host doesn't persist style-collection renames.
$wb.Styles.Item(1).Name = "Standard"

# Final selection left on the sheet
$ws.Range("F14").Select()
